$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.253.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.970.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.34%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -12.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.596"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "53.87"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.43%  "
$ws.Range("E9").Value = "  -4.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0745"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.00%  "
$ws.Range("E12").Value = "  -3.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.260.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "13.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.751"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.968.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.184.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.41%  "
$ws.Range("E20").Value = "  -3.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0801"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "220.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.32%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("E26").Value = "  -13.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.123"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("E32").Value = "  -3.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0600"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.75%  "
$ws.Range("E35").Value = "  -6.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.56%  "
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("E38").Value = "  -3.10%  "
$ws.Range("E39").Value = "  -4.21%  "
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.428.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.87%  "
$ws.Range("E43").Value = "  -6.60%  "
$ws.Range("E44").Value = "  -6.70%  "
$ws.Range("E45").Value = "  -11.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.983"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "14.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.28%  "
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("E50").Value = "  -4.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.50%  "
